# Natmi following Dr Hou advice
# Rebuild the Efna1 -> Epha7 LR-pair table for all 3x4 sending/target cluster
# combinations (ECs/FAPs/sCs sending; ECs/FAPs/M2/sCs target), replacing the
# previous partial (6-row) table with the full 12-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efna1"
$ws.Cells.Item(2, 3).Value = "Epha7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 23.60223533333333
$ws.Cells.Item(2, 8).Value = 70.80670599999999
$ws.Cells.Item(2, 9).Value = 0.8824726436021215
$ws.Cells.Item(2, 10).Value = 0.8824726436021214
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.03360366666666666
$ws.Cells.Item(2, 14).Value = 0.100811
$ws.Cells.Item(2, 15).Value = 0.0235667354422325
$ws.Cells.Item(2, 16).Value = 0.0235667354422325
$ws.Cells.Item(2, 17).Value = 0.7931216487295554
$ws.Cells.Item(2, 18).Value = 7.138094838565999
$ws.Cells.Item(2, 19).Value = 0.02079699932677872
$ws.Cells.Item(2, 20).Value = 0.02079699932677872

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efna1"
$ws.Cells.Item(3, 3).Value = "Epha7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 23.60223533333333
$ws.Cells.Item(3, 8).Value = 70.80670599999999
$ws.Cells.Item(3, 9).Value = 0.8824726436021215
$ws.Cells.Item(3, 10).Value = 0.8824726436021214
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 1.151356333333333
$ws.Cells.Item(3, 14).Value = 3.454069
$ws.Cells.Item(3, 15).Value = 0.8074627800757513
$ws.Cells.Item(3, 16).Value = 0.8074627800757513
$ws.Cells.Item(3, 17).Value = 27.17458313185711
$ws.Cells.Item(3, 18).Value = 244.571248186714
$ws.Cells.Item(3, 19).Value = 0.7125638141437668
$ws.Cells.Item(3, 20).Value = 0.7125638141437667

# Row 4: ECs -> M2
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efna1"
$ws.Cells.Item(4, 3).Value = "Epha7"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 23.60223533333333
$ws.Cells.Item(4, 8).Value = 70.80670599999999
$ws.Cells.Item(4, 9).Value = 0.8824726436021215
$ws.Cells.Item(4, 10).Value = 0.8824726436021214
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.009467333333333333
$ws.Cells.Item(4, 14).Value = 0.028402
$ws.Cells.Item(4, 15).Value = 0.00663957722897588
$ws.Cells.Item(4, 16).Value = 0.00663957722897588
$ws.Cells.Item(4, 17).Value = 0.2234502293124444
$ws.Cells.Item(4, 18).Value = 2.011052063812
$ws.Cells.Item(4, 19).Value = 0.005859245269654793
$ws.Cells.Item(4, 20).Value = 0.005859245269654792

# Row 5: ECs -> sCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Efna1"
$ws.Cells.Item(5, 3).Value = "Epha7"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 23.60223533333333
$ws.Cells.Item(5, 8).Value = 70.80670599999999
$ws.Cells.Item(5, 9).Value = 0.8824726436021215
$ws.Cells.Item(5, 10).Value = 0.8824726436021214
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 0.2314666666666667
$ws.Cells.Item(5, 14).Value = 0.6944
$ws.Cells.Item(5, 15).Value = 0.1623309072530403
$ws.Cells.Item(5, 16).Value = 0.1623309072530403
$ws.Cells.Item(5, 17).Value = 5.463130738488888
$ws.Cells.Item(5, 18).Value = 49.16817664639999
$ws.Cells.Item(5, 19).Value = 0.1432525848619213
$ws.Cells.Item(5, 20).Value = 0.1432525848619213

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efna1"
$ws.Cells.Item(6, 3).Value = "Epha7"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 2.633202
$ws.Cells.Item(6, 8).Value = 7.899606
$ws.Cells.Item(6, 9).Value = 0.09845375648791208
$ws.Cells.Item(6, 10).Value = 0.09845375648791205
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.03360366666666666
$ws.Cells.Item(6, 14).Value = 0.100811
$ws.Cells.Item(6, 15).Value = 0.0235667354422325
$ws.Cells.Item(6, 16).Value = 0.0235667354422325
$ws.Cells.Item(6, 17).Value = 0.088485242274
$ws.Cells.Item(6, 18).Value = 0.796367180466
$ws.Cells.Item(6, 19).Value = 0.002320233632444605
$ws.Cells.Item(6, 20).Value = 0.002320233632444604

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efna1"
$ws.Cells.Item(7, 3).Value = "Epha7"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 2.633202
$ws.Cells.Item(7, 8).Value = 7.899606
$ws.Cells.Item(7, 9).Value = 0.09845375648791208
$ws.Cells.Item(7, 10).Value = 0.09845375648791205
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 1.151356333333333
$ws.Cells.Item(7, 14).Value = 3.454069
$ws.Cells.Item(7, 15).Value = 0.8074627800757513
$ws.Cells.Item(7, 16).Value = 0.8074627800757513
$ws.Cells.Item(7, 17).Value = 3.031753799646001
$ws.Cells.Item(7, 18).Value = 27.28578419681401
$ws.Cells.Item(7, 19).Value = 0.07949774392263052
$ws.Cells.Item(7, 20).Value = 0.0794977439226305

# Row 8: FAPs -> M2
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Efna1"
$ws.Cells.Item(8, 3).Value = "Epha7"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 2.633202
$ws.Cells.Item(8, 8).Value = 7.899606
$ws.Cells.Item(8, 9).Value = 0.09845375648791208
$ws.Cells.Item(8, 10).Value = 0.09845375648791205
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009467333333333333
$ws.Cells.Item(8, 14).Value = 0.028402
$ws.Cells.Item(8, 15).Value = 0.00663957722897588
$ws.Cells.Item(8, 16).Value = 0.00663957722897588
$ws.Cells.Item(8, 17).Value = 0.024929401068
$ws.Cells.Item(8, 18).Value = 0.224364609612
$ws.Cells.Item(8, 19).Value = 0.0006536913196842773
$ws.Cells.Item(8, 20).Value = 0.0006536913196842771

# Row 9: FAPs -> sCs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Efna1"
$ws.Cells.Item(9, 3).Value = "Epha7"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 2.633202
$ws.Cells.Item(9, 8).Value = 7.899606
$ws.Cells.Item(9, 9).Value = 0.09845375648791208
$ws.Cells.Item(9, 10).Value = 0.09845375648791205
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 0.2314666666666667
$ws.Cells.Item(9, 14).Value = 0.6944
$ws.Cells.Item(9, 15).Value = 0.1623309072530403
$ws.Cells.Item(9, 16).Value = 0.1623309072530403
$ws.Cells.Item(9, 17).Value = 0.6094984896000001
$ws.Cells.Item(9, 18).Value = 5.485486406400001
$ws.Cells.Item(9, 19).Value = 0.01598208761315267
$ws.Cells.Item(9, 20).Value = 0.01598208761315267

# Row 10: sCs -> ECs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efna1"
$ws.Cells.Item(10, 3).Value = "Epha7"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.5101343333333334
$ws.Cells.Item(10, 8).Value = 1.530403
$ws.Cells.Item(10, 9).Value = 0.01907359990996641
$ws.Cells.Item(10, 10).Value = 0.0190735999099664
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.03360366666666666
$ws.Cells.Item(10, 14).Value = 0.100811
$ws.Cells.Item(10, 15).Value = 0.0235667354422325
$ws.Cells.Item(10, 16).Value = 0.0235667354422325
$ws.Cells.Item(10, 17).Value = 0.01714238409255555
$ws.Cells.Item(10, 18).Value = 0.154281456833
$ws.Cells.Item(10, 19).Value = 0.0004495024830091679
$ws.Cells.Item(10, 20).Value = 0.0004495024830091677

# Row 11: sCs -> FAPs
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Efna1"
$ws.Cells.Item(11, 3).Value = "Epha7"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 0.5101343333333334
$ws.Cells.Item(11, 8).Value = 1.530403
$ws.Cells.Item(11, 9).Value = 0.01907359990996641
$ws.Cells.Item(11, 10).Value = 0.0190735999099664
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 1.151356333333333
$ws.Cells.Item(11, 14).Value = 3.454069
$ws.Cells.Item(11, 15).Value = 0.8074627800757513
$ws.Cells.Item(11, 16).Value = 0.8074627800757513
$ws.Cells.Item(11, 17).Value = 0.5873463955341112
$ws.Cells.Item(11, 18).Value = 5.286117559807001
$ws.Cells.Item(11, 19).Value = 0.01540122200935408
$ws.Cells.Item(11, 20).Value = 0.01540122200935407

# Row 12: sCs -> M2
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Efna1"
$ws.Cells.Item(12, 3).Value = "Epha7"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 0.5101343333333334
$ws.Cells.Item(12, 8).Value = 1.530403
$ws.Cells.Item(12, 9).Value = 0.01907359990996641
$ws.Cells.Item(12, 10).Value = 0.0190735999099664
$ws.Cells.Item(12, 11).Value = 1.0
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.009467333333333333
$ws.Cells.Item(12, 14).Value = 0.028402
$ws.Cells.Item(12, 15).Value = 0.00663957722897588
$ws.Cells.Item(12, 16).Value = 0.00663957722897588
$ws.Cells.Item(12, 17).Value = 0.004829611778444444
$ws.Cells.Item(12, 18).Value = 0.043466506006
$ws.Cells.Item(12, 19).Value = 0.0001266406396368094
$ws.Cells.Item(12, 20).Value = 0.0001266406396368093

# Row 13: sCs -> sCs
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Efna1"
$ws.Cells.Item(13, 3).Value = "Epha7"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 0.5101343333333334
$ws.Cells.Item(13, 8).Value = 1.530403
$ws.Cells.Item(13, 9).Value = 0.01907359990996641
$ws.Cells.Item(13, 10).Value = 0.0190735999099664
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 0.2314666666666667
$ws.Cells.Item(13, 14).Value = 0.6944
$ws.Cells.Item(13, 15).Value = 0.1623309072530403
$ws.Cells.Item(13, 16).Value = 0.1623309072530403
$ws.Cells.Item(13, 17).Value = 0.1180790936888889
$ws.Cells.Item(13, 18).Value = 1.0627118432
$ws.Cells.Item(13, 19).Value = 0.003096234777966356
$ws.Cells.Item(13, 20).Value = 0.003096234777966354

